{"js": "// 1. Guppies paragraph: add a trailing \":\" run after \"(Martin and Johnsen 2007)\"\n{\n  const results = context.document.body.search(\"(Martin and Johnsen 2007)\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[0];\n  const para = found.paragraphs.getFirst();\n  para.insertText(\":\", Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 2. Barn swallows paragraph: \"** Barn swallows**\" -> bold \"Barn swallows\",\n//    and add a trailing \":\" run after \"(M\u00f8ller 1990)\"\n{\n  const results = context.document.body.search(\"** Barn swallows**\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[0];\n  const newRange = found.insertText(\"Barn swallows\", Word.InsertLocation.replace);\n  newRange.font.bold = true;\n  await context.sync();\n}\n{\n  const results = context.document.body.search(\"(M\u00f8ller 1990)\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[0];\n  const para = found.paragraphs.getFirst();\n  para.insertText(\":\", Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 3. Bibliography entry: italicize \"Poecilia Reticulata)\" within the Martin/Johnsen citation\n{\n  const results = context.document.body.search(\"Poecilia Reticulata)\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[0];\n  found.font.italic = true;\n  await context.sync();\n}\n\n// 4. Update the \"Last updated\" timestamp\n{\n  const results = context.document.body.search(\"Last updated: 2022-03-20 18:00:04\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  const found = results.items[0];\n  found.insertText(\"Last updated: 2022-03-20 21:11:47\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Guppies paragraph: append a trailing \":\" run right after \"(Martin and Johnsen 2007)\"\n$find1 = $d.Content.Find\n$find1.Text = \"(Martin and Johnsen 2007)\"\n$find1.Execute() | Out-Null\nif ($find1.Found) {\n    $find1.Parent.InsertAfter(\":\")\n}\n\n# 2. Barn swallows paragraph:\n#    \"** Barn swallows**\" -> bold \"Barn swallows\" (markdown asterisks removed),\n#    then append a trailing \":\" run after \"(M\u00f8ller 1990)\"\n$find2 = $d.Content.Find\n$find2.Text = \"** Barn swallows**\"\n$find2.Execute() | Out-Null\nif ($find2.Found) {\n    $r2 = $find2.Parent\n    $r2.Text = \"Barn swallows\"\n    $r2.Font.Bold = 1\n}\n\n$find3 = $d.Content.Find\n$find3.Text = \"(M\u00f8ller 1990)\"\n$find3.Execute() | Out-Null\nif ($find3.Found) {\n    $find3.Parent.InsertAfter(\":\")\n}\n\n# 3. Bibliography entry: italicize \"Poecilia Reticulata)\" within the Martin/Johnsen citation\n$find4 = $d.Content.Find\n$find4.Text = \"Poecilia Reticulata)\"\n$find4.Execute() | Out-Null\nif ($find4.Found) {\n    $find4.Parent.Font.Italic = 1\n}\n\n# 4. Update the \"Last updated\" timestamp\n$find5 = $d.Content.Find\n$find5.Text = \"Last updated: 2022-03-20 18:00:04\"\n$find5.Execute() | Out-Null\nif ($find5.Found) {\n    $find5.Parent.Text = \"Last updated: 2022-03-20 21:11:47\"\n}\n"}
